# issue #5: add legislator_id, name, date into dataframe
# Target worksheet is "股票" (stocks), which holds the securities table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# New header cells (row 1): date, legislator_name, legislator_id
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# New data cells (row 2). H2 holds a date-formatted string ("2012-02-29") that
# must stay literal text rather than being auto-converted to a date serial,
# so format the cell as Text first.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2012-02-29"
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").Value = "林岱樺"
$ws.Range("J2").Value = 904
